{"js": "// fix: update CLS nofo\n// Replace the eligibility-requirements sentence with the revised wording.\nconst oldText =\n  \". In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services. Eligible applicants may include public agencies and nonprofit agencies that provide other services, but eligible applicants must include legal services in their core services. Eligible applicants must also demonstrate a record of providing effective direct services to crime victims.\";\n\nconst newText =\n  \". In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services through the use of attorneys. If a victim service agency doesn\\u2019t currently focus on the provision of legal services, then it is not eligible. \";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# fix: update CLS nofo\n# Replace the eligibility-requirements sentence with the revised wording.\n$d = $word.ActiveDocument\n\n$oldText = \". In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services. Eligible applicants may include public agencies and nonprofit agencies that provide other services, but eligible applicants must include legal services in their core services. Eligible applicants must also demonstrate a record of providing effective direct services to crime victims.\"\n$newText = \". In brief, eligible applicants include public agencies and nonprofit organizations, whose primary mission is the provision of legal assistance services through the use of attorneys. If a victim service agency doesn\u2019t currently focus on the provision of legal services, then it is not eligible. \"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n$result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\nif (-not $result) {\n    throw \"Target sentence not found in document.\"\n}\n\nWrite-Output \"done\"\n"}
